$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update database: quarterly report publish dates shifted from
# --- 1401-10-27 to the newer 1402-01-28 release (columns I, J, M on the
# --- "تاریخ انتشار" row).
$ws.Range("I9").Value = "1402-01-28 (5)"
$ws.Range("J9").Value = "1402-01-28 (8)"
$ws.Range("M9").Value = "1402-01-28 (3)"

# --- Change read_price algorithm: capital ("سرمایه") is now correctly
# --- populated for the three columns that used to read 0, which in turn
# --- lets "سود هر سهم پس از کسر مالیات" (EPS after tax) be computed as a
# --- real number instead of falling back to the "-" placeholder.
$ws.Range("I26").Value = 2000000
$ws.Range("J26").Value = 2000000
$ws.Range("M26").Value = 2000000

$ws.Range("I25").Value = 5786
$ws.Range("J25").Value = 7282
$ws.Range("M25").Value = 8322
